$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values (column D) are plain decimal numbers (e.g. "1.001",
# "20.00"). Excel auto-coerces such strings to numeric cells on assignment,
# which would lose the original text formatting (trailing zeros, exact
# decimal form) used throughout this sheet. Force those specific cells to
# Text format first so the values are stored verbatim as text, matching the
# source data.

$ws.Range("D2").Value = "27.928.80"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.762.25"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.04"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4651"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3511"
$ws.Range("E8").Value = "  -2.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.44"
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.082"
$ws.Range("E11").Value = "  -2.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9997"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.53"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.995"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.156"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "1.763.79"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.36"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06425"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9995"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.82"
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.759"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "27.951.25"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("E24").Value = "  -1.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.152"
$ws.Range("E25").Value = "  +3.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.56"
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.00"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").Value = "1.960.93"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.167"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.15"
$ws.Range("E30").Value = "  -2.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.069"
$ws.Range("E31").Value = "  -2.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09306"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.550"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06070"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2062"
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.897"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6133"
$ws.Range("E40").Value = "  -3.11%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.787"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.353"
$ws.Range("E43").Value = "  -3.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.11"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.732"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5787"
$ws.Range("E46").Value = "  -1.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.92"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.925"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06809"
$ws.Range("E49").Value = "  -2.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.122"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.14"
$ws.Range("E51").Value = "  -0.38%  "
